# Update countries & provincias Spain
# This script applies the daily COVID-19 data refresh to the "Pais" sheet:
#  - Updates the "Datos actualizados..." timestamp in A1
#  - Refreshes several countries' statistics (Estados Unidos, Alemania, Barein,
#    Uzbekistan, Bulgaria, Guinea Ecuatorial, etc.)
#  - Sudafrica's total cases overtook Panama/Australia/Egipto, so it moves up
#    the (descending, sorted-by-total-cases) table; the rows that used to hold
#    Panama/Australia/Egipto now shift down one position and keep their own
#    (unchanged) figures.
#  - San Cristobal y Nieves overtook Burundi, so those two adjacent rows swap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / timestamp -----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 21:03"

# --- Estados Unidos (row 4) --------------------------------------------------
$ws.Range("B4").Value = 1201460
$ws.Range("C4").Value = 13338
$ws.Range("D4").Value = 181603
$ws.Range("E4").Value = 950714
$ws.Range("G4").Value = 545
$ws.Range("H4").Value = 69143

# --- Alemania (row 9) ---------------------------------------------------------
$ws.Range("B9").Value = 165786
$ws.Range("C9").Value = 122
$ws.Range("E9").Value = 26193
$ws.Range("G9").Value = 27
$ws.Range("H9").Value = 6893

# --- Sudafrica moves up above Panama/Australia/Egipto (rows 49-52) ---------
# New Sudafrica row (inserted right after Colombia, pushing the rest down)
$ws.Range("A49").Value = "Sudafrica"
$ws.Range("B49").Value = 7220
$ws.Range("C49").Value = 437
$ws.Range("D49").Value = 2746
$ws.Range("E49").Value = 4336
$ws.Range("F49").Value = 36
$ws.Range("G49").Value = 7
$ws.Range("H49").Value = 138

# Panama (was row 49, now row 50) keeps its previous data
$ws.Range("A50").Value = "Panama"
$ws.Range("B50").Value = 7197
$ws.Range("C50").Value = 107
$ws.Range("D50").Value = 641
$ws.Range("E50").Value = 6356
$ws.Range("F50").Value = 91
$ws.Range("G50").Value = 3
$ws.Range("H50").Value = 200

# Australia (was row 50, now row 51) keeps its previous data
$ws.Range("A51").Value = "Australia"
$ws.Range("B51").Value = 6825
$ws.Range("C51").Value = 24
$ws.Range("D51").Value = 5859
$ws.Range("E51").Value = 871
$ws.Range("F51").Value = 28
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 95

# Egipto (was row 51, now row 52) keeps its previous data
$ws.Range("A52").Value = "Egipto"
$ws.Range("B52").Value = 6813
$ws.Range("C52").Value = 348
$ws.Range("D52").Value = 1632
$ws.Range("E52").Value = 4745
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 7
$ws.Range("H52").Value = 436

# --- Barein (row 62) ----------------------------------------------------------
$ws.Range("D62").Value = 1744
$ws.Range("E62").Value = 1781

# --- Uzbekistan (row 72) ------------------------------------------------------
$ws.Range("B72").Value = 2189
$ws.Range("C72").Value = 40
$ws.Range("E72").Value = 774

# --- Bulgaria (row 81) ---------------------------------------------------------
$ws.Range("F81").Value = 37

# --- Guinea Ecuatorial (row 130) ---------------------------------------------
$ws.Range("E130").Value = 299
$ws.Range("G130").Value = 2
$ws.Range("H130").Value = 3

# --- San Cristobal y Nieves overtakes Burundi (rows 198-199 swap) ----------
$ws.Range("A198").Value = "San Cristobal y Nieves"
$ws.Range("B198").Value = 15
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 8
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Burundi"
$ws.Range("B199").Value = 15
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 7
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 1
